# feat: add 2022-Q4 data
#
# Before:  "总计" (summary) + "2020-Q4" (fund holdings detail)
# After:   "总计" (summary, +1 row) + "2022-Q4" (NEW fund holdings detail)
#          + "2020-Q4" (unchanged fund holdings detail, moved to a new sheet part)
#
# Strategy: duplicate the existing "2020-Q4" sheet so the duplicate keeps the
# "2020-Q4" name/content untouched (placed after, becomes the new part), while
# the original sheet (in its original slot) is renamed to "2022-Q4" and has its
# data replaced with the new quarter's numbers. This reproduces the sheetId/
# r:id layout the diff shows (2022-Q4 keeps the old 2020-Q4 slot's id, the
# untouched 2020-Q4 data lands in a freshly minted slot).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value as literal text (no "smart" numeric/leading-zero
# coercion), then strip the NumberFormat override back to General so the
# saved cell carries no style index -- matching cells authored as plain
# inline strings.
# ---------------------------------------------------------------------------
function Set-TextValue($range, $text, $blankFormatCell) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $blankFormatCell.Copy()
    $range.PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# Step 1: duplicate "2020-Q4" (sheet index 2) so the copy (placed right after
# it) becomes the sheet that keeps the "2020-Q4" name and data.
# ---------------------------------------------------------------------------
$orig2020 = $wb.Worksheets.Item(2)
$orig2020.Copy($null, $orig2020)
$copy2020 = $wb.Worksheets.Item(3)

$orig2020.Name = "2022-Q4"
$copy2020.Name = "2020-Q4"

$ws22 = $orig2020
$wsTotal = $wb.Worksheets.Item(1)
$blank = $wsTotal.Range("C2")

# ---------------------------------------------------------------------------
# Step 2: drop the now-stale 5th data row from the duplicated sheet (the new
# 2022-Q4 sheet only has 3 data rows) before rewriting its contents.
# ---------------------------------------------------------------------------
$ws22.Rows.Item(5).Delete()

# ---------------------------------------------------------------------------
# Step 3: headers - D1 text changes ("基金金额" -> "基金规模"); all header
# cells (B1:H1) get the bordered/centered header style used elsewhere in this
# workbook (style copied cross-sheet from the "总计" header row).
# ---------------------------------------------------------------------------
$ws22.Range("D1").Value = "基金规模"

$wsTotal.Range("B1").Copy()
$ws22.Range("B1:H1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Step 4: replace the 3 data rows with the 2022-Q4 holdings.
# ---------------------------------------------------------------------------
$wsTotal.Range("A2").Copy()
$ws22.Range("A2:A4").PasteSpecial(-4122)

$ws22.Range("A2").Value = 0
Set-TextValue $ws22.Range("B2") "013442" $blank
$ws22.Range("C2").Value = "建信中证1000指数增强E"
Set-TextValue $ws22.Range("D2") "9.52" $blank
Set-TextValue $ws22.Range("E2") "86.80" $blank
Set-TextValue $ws22.Range("F2") "1.23" $blank
Set-TextValue $ws22.Range("G2") "0.1171" $blank
$ws22.Range("H2").Value = 7

$ws22.Range("A3").Value = 1
Set-TextValue $ws22.Range("B3") "006165" $blank
$ws22.Range("C3").Value = "建信中证1000指数增强A"
Set-TextValue $ws22.Range("D3") "7.20" $blank
Set-TextValue $ws22.Range("E3") "86.80" $blank
Set-TextValue $ws22.Range("F3") "1.23" $blank
Set-TextValue $ws22.Range("G3") "0.0886" $blank
$ws22.Range("H3").Value = 7

$ws22.Range("A4").Value = 2
Set-TextValue $ws22.Range("B4") "006166" $blank
$ws22.Range("C4").Value = "建信中证1000指数增强C"
Set-TextValue $ws22.Range("D4") "2.21" $blank
Set-TextValue $ws22.Range("E4") "86.80" $blank
Set-TextValue $ws22.Range("F4") "1.23" $blank
Set-TextValue $ws22.Range("G4") "0.0272" $blank
$ws22.Range("H4").Value = 7

# ---------------------------------------------------------------------------
# Step 5: insert the 2022-Q4 summary row on "总计" (row 2), pushing the
# existing 2020-Q4 row down to row 3.
# ---------------------------------------------------------------------------
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2020-Q4"
$wsTotal.Range("C3").Value = 4
$wsTotal.Range("D3").Value = 0.89

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 3
$wsTotal.Range("D2").Value = 0.23
